# MesterLista.xlsx - "Lista" sheet edit
#
# Summary of the target change (see commit message: "new file: barcode
# generator for labels"):
#   - The City ("Varos"/B) and Country ("Orszag"/C) columns are swapped
#     (and likewise their English counterparts, "Angol Varos"/E and
#     "Angol Orszag"/F) for the header row and all 100 data rows.
#   - A handful of incidental layout changes (slightly narrower A/D and
#     wider B/C, E/F columns, a larger used range) accompany the data
#     reshuffle; we reproduce what is reachable from the Excel object
#     model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 101

for ($r = 1; $r -le $lastRow; $r++) {
    $city    = $ws.Cells.Item($r, 2).Value2   # B: Varos
    $country = $ws.Cells.Item($r, 3).Value2   # C: Orszag
    $ws.Cells.Item($r, 2).Value = $country
    $ws.Cells.Item($r, 3).Value = $city

    $cityEn    = $ws.Cells.Item($r, 5).Value2 # E: Angol Varos
    $countryEn = $ws.Cells.Item($r, 6).Value2 # F: Angol Orszag
    $ws.Cells.Item($r, 5).Value = $countryEn
    $ws.Cells.Item($r, 6).Value = $cityEn
}

# Nudge the sheet's used range out to column Q (matches the new
# <dimension ref="A1:Q101"/>) without disturbing any cell content.
$ws.Range("Q101").Font.Bold = $false

# Best-effort column width tweaks to track the target layout.
$ws.Columns.Item(1).ColumnWidth = 38
$ws.Columns.Item(2).ColumnWidth = 23.8333333333333
$ws.Columns.Item(3).ColumnWidth = 20.6666666666667
$ws.Columns.Item(4).ColumnWidth = 38
$ws.Columns.Item(5).ColumnWidth = 23.8333333333333
$ws.Columns.Item(6).ColumnWidth = 20.6666666666667
$ws.Columns.Item(7).ColumnWidth = 10.3333333333333
$ws.Columns.Item(8).ColumnWidth = 10.3333333333333
$ws.Columns.Item(9).ColumnWidth = 10.3333333333333
